$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25)
$Bvals = @(15.6333735743518, 14.86797853438551, 14.37878427975039, 14.17483628783648, 14.14070100983407, 14.37605203395405, 15.37359075367866, 17.16839700975018, 18.37921429312474, 18.90516815627925, 19.10066536695179, 19.05872611873617, 18.92132587215132, 18.83668374550258, 18.34433322051014, 18.0358523061863, 17.85608745634046, 17.79482439881592, 18.06893302232145, 18.9617839399954, 19.5238887899039, 19.22587056246032, 18.05398474546012, 16.70121519621702)
$Cvals = @(10.77385563564584, 10.12104376140404, 9.6964385397083, 9.517461930409029, 9.487385002728882, 9.694048825362851, 10.5537168695809, 12.04995978898853, 13.03302162647491, 13.45498127092831, 13.61113936901156, 13.57766927926213, 13.46790120443458, 13.40019257806018, 13.00493765466055, 12.7559963433615, 12.61043589103537, 12.56074357907673, 12.78274248280629, 13.50024121417377, 13.94801875218763, 13.71096561264703, 12.77065815575815, 11.66557346941336)
$Dvals = @(5.968165594392101, 5.84689993182348, 5.773059609588336, 5.743170413395644, 5.738220849999554, 5.772655637005927, 5.926252590775858, 6.230360725910545, 6.452970670417021, 6.553547596103231, 6.591494027042382, 6.583328428779922, 6.556672500060177, 6.540325589638753, 6.446380324919181, 6.388541669082916, 6.35521180631063, 6.343917414299844, 6.394705473583352, 6.564506101597837, 6.674648130764202, 6.6159522840645, 6.391919056334736, 6.148053482900536)
$Evals = @(16.53126690945976, 15.58853288814663, 14.98503997366339, 14.73318237601567, 14.69101241570584, 14.98166694662424, 16.2114741812632, 18.51775793847523, 20.17319181352013, 20.88430613434355, 21.14759292288251, 21.09115554363144, 20.90608687723363, 20.79194707529607, 20.1258777188151, 19.7065465927658, 19.46140768340027, 19.3777282272498, 19.75159373286345, 20.96060849769838, 21.71583220212232, 21.3159401168671, 19.73124056960713, 17.87069104989532)
$Gvals = @(3.606950018463854, 3.610008915970679, 3.611981970187542, 3.612809956495366, 3.612948892272536, 3.611993039603339, 3.607985090168384, 3.600874095281281, 3.596100052599133, 3.594024743711961, 3.593252646344426, 3.593418320003365, 3.593960947193427, 3.594295113556885, 3.596237611659424, 3.597453905195138, 3.59816256655066, 3.598404069577771, 3.597323489423087, 3.593801191237288, 3.591579427890479, 3.592757910168205, 3.597382421112315, 3.602718283000311)
$Ovals = @(20.38903314766108, 20.41558900619489, 20.44100398404924, 20.45363317324225, 20.45586691406482, 20.44116513072634, 20.3962885068153, 20.38128681596864, 20.41560622515734, 20.44119522310502, 20.45232667034793, 20.44986512417275, 20.44208209537013, 20.43750264266128, 20.41413532243356, 20.40236118249879, 20.39652783803772, 20.39471376463754, 20.40351733772573, 20.44432899760748, 20.47940499351759, 20.45991379604136, 20.40299172561344, 20.37743450856214)

for ($i = 0; $i -lt $rows.Length; $i++) {
  $r = $rows[$i]
  $ws.Cells.Item($r, 2).Value = $Bvals[$i]
  $ws.Cells.Item($r, 3).Value = $Cvals[$i]
  $ws.Cells.Item($r, 4).Value = $Dvals[$i]
  $ws.Cells.Item($r, 5).Value = $Evals[$i]
  $ws.Cells.Item($r, 7).Value = $Gvals[$i]
  $ws.Cells.Item($r, 15).Value = $Ovals[$i]
}
